$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (prevents Excel from
# re-interpreting dotted numeric-looking strings like "604.55"
# or "63.726.23" as numbers/dates), then restore the default style
# so no stray cell formatting is introduced.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '63.726.23'
$ws.Range('E2').Value = '  +1.27%  '
Set-TextValue $ws.Range('D3') '3.317.71'
$ws.Range('E3').Value = '  +4.67%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue $ws.Range('D5') '604.55'
$ws.Range('E5').Value = '  +2.71%  '
Set-TextValue $ws.Range('D6') '142.41'
$ws.Range('E6').Value = '  +2.80%  '
$ws.Range('E7').Value = '  +0.05%  '
Set-TextValue $ws.Range('D8') '3.316.45'
$ws.Range('E8').Value = '  +4.64%  '
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('E13').Value = '  +1.28%  '
Set-TextValue $ws.Range('D14') '34.87'
$ws.Range('E14').Value = '  +2.61%  '
Set-TextValue $ws.Range('D15') '3.863.25'
$ws.Range('E15').Value = '  +4.76%  '
Set-TextValue $ws.Range('D16') '0.120'
Set-TextValue $ws.Range('D17') '3.313.51'
$ws.Range('E17').Value = '  +4.63%  '
Set-TextValue $ws.Range('D18') '63.799.29'
$ws.Range('E19').Value = '  +3.08%  '
Set-TextValue $ws.Range('D20') '480.22'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('E22').Value = '  +4.87%  '
Set-TextValue $ws.Range('D23') '7.95'
$ws.Range('E23').Value = '  +2.87%  '
Set-TextValue $ws.Range('D24') '13.72'
$ws.Range('E24').Value = '  +5.66%  '
Set-TextValue $ws.Range('D25') '84.75'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  -0.01%  '
Set-TextValue $ws.Range('D27') '2.78'
$ws.Range('E27').Value = '  +2.70%  '
$ws.Range('E28').Value = '  -0.04%  '
Set-TextValue $ws.Range('D29') '8.23'
$ws.Range('E29').Value = '  +3.56%  '
Set-TextValue $ws.Range('D30') '7.19'
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('E31').Value = '  +3.92%  '
Set-TextValue $ws.Range('D32') '28.93'
$ws.Range('E32').Value = '  +7.93%  '
Set-TextValue $ws.Range('D33') '0.106'
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('E35').Value = '  +2.74%  '
Set-TextValue $ws.Range('D36') '6.09'
$ws.Range('E36').Value = '  +5.43%  '
Set-TextValue $ws.Range('D37') '52.37'
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('E38').Value = '  +5.16%  '
$ws.Range('E39').Value = '  +3.57%  '
Set-TextValue $ws.Range('D40') '435.65'
$ws.Range('E40').Value = '  +4.70%  '
Set-TextValue $ws.Range('D41') '3.101.75'
$ws.Range('E41').Value = '  +4.91%  '
Set-TextValue $ws.Range('D42') '0.119'
$ws.Range('E42').Value = '  +7.79%  '
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('E45').Value = '  +1.63%  '
$ws.Range('E46').Value = '  +5.87%  '
Set-TextValue $ws.Range('D47') '37.04'
$ws.Range('E47').Value = '  +15.72%  '
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  +2.64%  '
$ws.Range('E51').Value = '  +0.25%  '
